$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the Shipment controller's integration-test column ("Yes") for all of
# its endpoints (rows 74-78: GET/POST/GET{id}/PUT{id}/DELETE{id} shipments)
$ws.Range("C74:C78").Value = "Yes"

# Leave the selection where the author left it when saving the workbook
$ws.Range("C78").Select()
